$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C2:C9): 45183 -> 45184
$ws.Range("C2:C9").Value = 45184

# Row 2 ("A 31987-2023") hyperlink formulas gain a visible-text second
# argument. S2 reproduces the source data's broken quoting verbatim
# (closing quote of the URL literal got moved to the very end).
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/artfynd/A 31987-2023.xlsx, "A 31987-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/kartor/A 31987-2023.png", "A 31987-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/klagomål/A 31987-2023.docx", "A 31987-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/klagomålsmail/A 31987-2023.docx", "A 31987-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/tillsyn/A 31987-2023.docx", "A 31987-2023")'

# Y2 used to be a plain inline string holding the literal (unexecuted)
# formula text with a semicolon argument separator; it becomes a real
# formula using a comma separator, like the other link columns.
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/tillsynsmail/A 31987-2023.docx", "A 31987-2023")'
